# edit.ps1
# Commit: "update scripts wuth new tpm"
#
# The underlying NATMI TPM recomputation changed every numeric statistic in the
# LR-pairs table and introduced a new sending/target cluster, "Resolving-Mac"
# (previously only a target cluster; now also a sending cluster), which adds
# five new rows (17-21) to the Sheet1 grid and grows the used range from
# A1:T16 to A1:T21. Rather than patch individual cells we rewrite the full
# data block (rows 2-21, columns A:T) with the post-recalculation values taken
# from the target workbook, which is simplest + least error-prone given how
# thoroughly every row changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner tuple is one full row: A..T in order. Columns A-D are text
# (cluster/gene names), columns E-T are numeric statistics.
$data = @(
    ,("ECs", "Myl9", "Cd69", "ECs", "3", "1", "5.009748000000001", "15.029244", "0.1179910219624604", "0.1179910219624604", "2", "0.6666666666666666", "0.1205896666666667", "0.361769", "0.003738558763786432", "0.003738558763786432", "0.6041238414040001", "5.437114572636001", "0.0004411163692058737", "0.0004411163692058738")
    ,("ECs", "Myl9", "Cd69", "FAPs", "3", "1", "5.009748000000001", "15.029244", "0.1179910219624604", "0.1179910219624604", "2", "0.6666666666666666", "0.05518799999999999", "0.165564", "0.001710955729118683", "0.001710955729118683", "0.276477972624", "2.488301753616", "0.00020187741501124", "0.00020187741501124")
    ,("ECs", "Myl9", "Cd69", "Inflammatory-Mac", "3", "1", "5.009748000000001", "15.029244", "0.1179910219624604", "0.1179910219624604", "3", "1", "22.44608266666667", "67.33824799999999", "0.6958805126985015", "0.6958805126985015", "112.449217747168", "1012.042959724512", "0.0821076528570571", "0.0821076528570571")
    ,("ECs", "Myl9", "Cd69", "MuSCs", "3", "1", "5.009748000000001", "15.029244", "0.1179910219624604", "0.1179910219624604", "3", "1", "0.104803", "0.314409", "0.003249135559883043", "0.003249135559883043", "0.5250366196440002", "4.725329576796002", "0.0003833688252051712", "0.0003833688252051713")
    ,("ECs", "Myl9", "Cd69", "Resolving-Mac", "3", "1", "5.009748000000001", "15.029244", "0.1179910219624604", "0.1179910219624604", "3", "1", "9.528993", "28.586979", "0.2954208372487103", "0.2954208372487104", "47.73785362376401", "429.640682613876", "0.03485700649598102", "0.03485700649598103")
    ,("FAPs", "Myl9", "Cd69", "ECs", "3", "1", "26.708345", "80.125035", "0.6290426028367033", "0.6290426028367033", "2", "0.6666666666666666", "0.1205896666666667", "0.361769", "0.003738558763786432", "0.003738558763786432", "3.220750420768333", "28.986753786915", "0.002351712735630185", "0.002351712735630185")
    ,("FAPs", "Myl9", "Cd69", "FAPs", "3", "1", "26.708345", "80.125035", "0.6290426028367033", "0.6290426028367033", "2", "0.6666666666666666", "0.05518799999999999", "0.165564", "0.001710955729118683", "0.001710955729118683", "1.47398014386", "13.26582129474", "0.001076264045183186", "0.001076264045183186")
    ,("FAPs", "Myl9", "Cd69", "Inflammatory-Mac", "3", "1", "26.708345", "80.125035", "0.6290426028367033", "0.6290426028367033", "3", "1", "22.44608266666667", "67.33824799999999", "0.6958805126985015", "0.6958805126985015", "599.4977197598532", "5395.479477838679", "0.4377384889712049", "0.4377384889712049")
    ,("FAPs", "Myl9", "Cd69", "MuSCs", "3", "1", "26.708345", "80.125035", "0.6290426028367033", "0.6290426028367033", "3", "1", "0.104803", "0.314409", "0.003249135559883043", "0.003249135559883043", "2.799114681035", "25.192032129315", "0.002043844689558118", "0.002043844689558119")
    ,("FAPs", "Myl9", "Cd69", "Resolving-Mac", "3", "1", "26.708345", "80.125035", "0.6290426028367033", "0.6290426028367033", "3", "1", "9.528993", "28.586979", "0.2954208372487103", "0.2954208372487104", "254.503632546585", "2290.532692919265", "0.1858322923951269", "0.1858322923951269")
    ,("MuSCs", "Myl9", "Cd69", "ECs", "3", "1", "10.339867", "31.019601", "0.2435275136166369", "0.2435275136166369", "2", "0.6666666666666666", "0.1205896666666667", "0.361769", "0.003738558763786432", "0.003738558763786432", "1.246881114907667", "11.221930034169", "0.0009104419202545975", "0.0009104419202545976")
    ,("MuSCs", "Myl9", "Cd69", "FAPs", "3", "1", "10.339867", "31.019601", "0.2435275136166369", "0.2435275136166369", "2", "0.6666666666666666", "0.05518799999999999", "0.165564", "0.001710955729118683", "0.001710955729118683", "0.5706365799959999", "5.135729219963999", "0.0004166647946204129", "0.000416664794620413")
    ,("MuSCs", "Myl9", "Cd69", "Inflammatory-Mac", "3", "1", "10.339867", "31.019601", "0.2435275136166369", "0.2435275136166369", "3", "1", "22.44608266666667", "67.33824799999999", "0.6958805126985015", "0.6958805126985015", "232.0895094443387", "2088.805584999048", "0.1694660510317366", "0.1694660510317366")
    ,("MuSCs", "Myl9", "Cd69", "MuSCs", "3", "1", "10.339867", "31.019601", "0.2435275136166369", "0.2435275136166369", "3", "1", "0.104803", "0.314409", "0.003249135559883043", "0.003249135559883043", "1.083649081201", "9.752841730809001", "0.0007912539043017169", "0.0007912539043017170")
    ,("MuSCs", "Myl9", "Cd69", "Resolving-Mac", "3", "1", "10.339867", "31.019601", "0.2435275136166369", "0.2435275136166369", "3", "1", "9.528993", "28.586979", "0.2954208372487103", "0.2954208372487104", "98.528520263931", "886.7566823753789", "0.07194310196572357", "0.07194310196572358")
    ,("Resolving-Mac", "Myl9", "Cd69", "ECs", "2", "0.6666666666666666", "0.400762", "1.202286", "0.009438861584199355", "0.009438861584199355", "2", "0.6666666666666666", "0.1205896666666667", "0.361769", "0.003738558763786432", "0.003738558763786432", "0.04832775599266666", "0.434949803934", "3.528773869577559E-05", "3.528773869577559E-05")
    ,("Resolving-Mac", "Myl9", "Cd69", "FAPs", "2", "0.6666666666666666", "0.400762", "1.202286", "0.009438861584199355", "0.009438861584199355", "2", "0.6666666666666666", "0.05518799999999999", "0.165564", "0.001710955729118683", "0.001710955729118683", "0.022117253256", "0.199055279304", "1.614947430384413E-05", "1.614947430384414E-05")
    ,("Resolving-Mac", "Myl9", "Cd69", "Inflammatory-Mac", "2", "0.6666666666666666", "0.400762", "1.202286", "0.009438861584199355", "0.009438861584199355", "3", "1", "22.44608266666667", "67.33824799999999", "0.6958805126985015", "0.6958805126985015", "8.995536981658667", "80.95983283492799", "0.006568319838502837", "0.006568319838502837")
    ,("Resolving-Mac", "Myl9", "Cd69", "MuSCs", "2", "0.6666666666666666", "0.400762", "1.202286", "0.009438861584199355", "0.009438861584199355", "3", "1", "0.104803", "0.314409", "0.003249135559883043", "0.003249135559883043", "0.04200105988600001", "0.378009538974", "3.066814081803612E-05", "3.066814081803612E-05")
    ,("Resolving-Mac", "Myl9", "Cd69", "Resolving-Mac", "2", "0.6666666666666666", "0.400762", "1.202286", "0.009438861584199355", "0.009438861584199355", "3", "1", "9.528993", "28.586979", "0.2954208372487103", "0.2954208372487104", "3.818858292666", "34.369724633994", "0.002788436391878862", "0.002788436391878862")
)

$startRow = 2
$numStringCols = 4   # columns A-D are text

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $col = $c + 1
        $cellValue = $rowValues[$c]
        if ($c -lt $numStringCols) {
            $ws.Cells.Item($row, $col).Value = $cellValue
        } else {
            $ws.Cells.Item($row, $col).Value = [double]$cellValue
        }
    }
}

